$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "44.221.07"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  +1.44%  "
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.254.88"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  +0.72%  "
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "308.19"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -4.33%  "
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "99.05"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -2.08%  "
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -0.51%  "
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  +0.10%  "
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -2.85%  "
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "35.69"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -3.44%  "
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0824"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -0.68%  "
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "7.34"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -4.57%  "
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -1.68%  "
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.598.16"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  +0.65%  "
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.256.13"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  +0.59%  "
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -1.46%  "
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "13.87"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -1.59%  "
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "44.081.45"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  +1.26%  "
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.08"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -3.09%  "
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0979"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -0.52%  "
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.34"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -3.09%  "
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "65.45"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  +0.78%  "
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "242.93"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  +2.71%  "
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.96"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -6.68%  "
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.99"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -7.47%  "
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  +0.21%  "
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.14"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  +1.10%  "
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.15"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -1.73%  "
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "36.63"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  +0.35%  "
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -1.84%  "
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "20.21"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  +0.36%  "
$cell.Style = "Normal"
$cell = $ws.Range("B32")
$cell.NumberFormat = "@"
$cell.Value = "LidoDAOToken"
$cell.Style = "Normal"
$cell = $ws.Range("C32")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.58"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +14.67%  "
$cell.Style = "Normal"
$cell = $ws.Range("B33")
$cell.NumberFormat = "@"
$cell.Value = "Monero"
$cell.Style = "Normal"
$cell = $ws.Range("C33")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "156.63"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -1.72%  "
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0829"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -2.31%  "
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  -0.96%  "
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  +0.15%  "
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -4.34%  "
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.87"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -2.34%  "
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "15.34"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -0.33%  "
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -7.82%  "
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "3.39"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -9.71%  "
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.0308"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -2.85%  "
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.776.33"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -0.58%  "
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "88.37"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  +7.41%  "
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  +13.79%  "
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +0.15%  "
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.194"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -2.71%  "
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "101.85"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -0.84%  "
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "8.29"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -1.28%  "
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "70.48"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -5.19%  "
$cell.Style = "Normal"
